$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 4301
$ws.Range("R2").Value = 0
$ws.Range("Y2").Value = "206"
$ws.Range("Z2").Value = "4506"
$ws.Range("AA2").Value = "406 x 140 x 53"
$ws.Range("AC2").Value = 4821
$ws.Range("AD2").Value = 3024
$ws.Range("AE2").Value = 1.26
$ws.Range("AG2").Value = 0.49
$ws.Range("AH2").Value = 2371
$ws.Range("AJ2").Value = 4821
$ws.Range("AK2").Value = 339500
$ws.Range("AL2").Value = 65057
$ws.Range("AP2").Value = 4743
$ws.Range("AR2").Value = 2410
$ws.Range("AS2").Value = 20564
$ws.Range("AT2").Value = 0.34
$ws.Range("AV2").Value = 0.97
$ws.Range("AW2").Value = 2332
$ws.Range("AX2").Value = 0.98
$ws.Range("AZ2").Value = 56095
$ws.Range("BA2").Value = 3
$ws.Range("BB2").Value = 25
$ws.Range("BC2").Value = 47
$ws.Range("BD2").Value = 357
$ws.Range("BE2").Value = 2289

# Row 3
$ws.Range("Q3").Value = 4562
$ws.Range("R3").Value = 0
$ws.Range("Y3").Value = "208"
$ws.Range("Z3").Value = "4770"
$ws.Range("AA3").Value = "406 x 178 x 54"
$ws.Range("AC3").Value = 4899
$ws.Range("AD3").Value = 12359
$ws.Range("AE3").Value = 0.63
$ws.Range("AG3").Value = 0.88
$ws.Range("AH3").Value = 4303
$ws.Range("AJ3").Value = 4899
$ws.Range("AK3").Value = 345000
$ws.Range("AL3").Value = 228014
$ws.Range("AP3").Value = 4956
$ws.Range("AR3").Value = 2450
$ws.Range("AS3").Value = 33032
$ws.Range("AW3").Value = 2410
$ws.Range("AX3").Value = 0.98
$ws.Range("AZ3").Value = 228014
$ws.Range("BA3").Value = 3
$ws.Range("BB3").Value = 13
$ws.Range("BC3").Value = 12
$ws.Range("BD3").Value = 152
$ws.Range("BE3").Value = 2400

# Row 4
$ws.Range("Q4").Value = 7095
$ws.Range("R4").Value = 0
$ws.Range("Z4").Value = "7095"
$ws.Range("AA4").Value = "533 x 210 x 82"
$ws.Range("AB4").Value = 355
$ws.Range("AC4").Value = 7455
$ws.Range("AD4").Value = 7848
$ws.Range("AE4").Value = 0.97
$ws.Range("AG4").Value = 0.68
$ws.Range("AH4").Value = 5094
$ws.Range("AI4").Value = 355
$ws.Range("AJ4").Value = 7455
$ws.Range("AK4").Value = 525000
$ws.Range("AL4").Value = 100603
$ws.Range("AP4").Value = 7335
$ws.Range("AQ4").Value = 355
$ws.Range("AR4").Value = 3728
$ws.Range("AS4").Value = 65093
$ws.Range("AT4").Value = 0.24
$ws.Range("AV4").Value = 0.99
$ws.Range("AW4").Value = 3695
$ws.Range("AX4").Value = 1.26
$ws.Range("AZ4").Value = 86745
$ws.Range("BA4").Value = 3
$ws.Range("BB4").Value = 25
$ws.Range("BC4").Value = 47
$ws.Range("BD4").Value = 531
$ws.Range("BE4").Value = 3601

# Row 5
$ws.Range("Q5").Value = 7525
$ws.Range("R5").Value = 0
$ws.Range("Z5").Value = "7525"
$ws.Range("AA5").Value = "533 x 210 x 92"
$ws.Range("AC5").Value = 8307
$ws.Range("AD5").Value = 36482
$ws.Range("AE5").Value = 0.48
$ws.Range("AG5").Value = 0.93
$ws.Range("AH5").Value = 7736
$ws.Range("AJ5").Value = 8307
$ws.Range("AK5").Value = 585000
$ws.Range("AL5").Value = 386633
$ws.Range("AP5").Value = 8403
$ws.Range("AR5").Value = 4154
$ws.Range("AS5").Value = 77399
$ws.Range("AT5").Value = 0.23
$ws.Range("AV5").Value = 0.99
$ws.Range("AW5").Value = 4125
$ws.Range("AX5").Value = 1.36
$ws.Range("AZ5").Value = 386633
$ws.Range("BA5").Value = 3
$ws.Range("BB5").Value = 13
$ws.Range("BC5").Value = 12
$ws.Range("BD5").Value = 232
$ws.Range("BE5").Value = 3786

# Row 6
$ws.Range("Q6").Value = 10164
$ws.Range("R6").Value = 0
$ws.Range("Z6").Value = "10164"
$ws.Range("AA6").Value = "610 x 229 x 125"
$ws.Range("AC6").Value = 10971
$ws.Range("AD6").Value = 16291
$ws.Range("AE6").Value = 0.82
$ws.Range("AG6").Value = 0.78
$ws.Range("AH6").Value = 8598
$ws.Range("AJ6").Value = 10971
$ws.Range("AK6").Value = 795000
$ws.Range("AL6").Value = 152342
$ws.Range("AP6").Value = 10804
$ws.Range("AR6").Value = 5486
$ws.Range("AS6").Value = 127272
$ws.Range("AT6").Value = 0.21
$ws.Range("AV6").Value = 1
$ws.Range("AW6").Value = 5476
$ws.Range("AX6").Value = 1.68
$ws.Range("AZ6").Value = 131356
$ws.Range("BA6").Value = 4
$ws.Range("BB6").Value = 34
$ws.Range("BC6").Value = 63
$ws.Range("BD6").Value = 748
$ws.Range("BE6").Value = 5157

# Row 7
$ws.Range("Q7").Value = 10781
$ws.Range("R7").Value = 0
$ws.Range("Z7").Value = "10781"
$ws.Range("AA7").Value = "686 x 254 x 125"
$ws.Range("AB7").Value = 345
$ws.Range("AC7").Value = 10971
$ws.Range("AD7").Value = 77987
$ws.Range("AE7").Value = 0.38
$ws.Range("AG7").Value = 0.96
$ws.Range("AH7").Value = 10524
$ws.Range("AI7").Value = 345
$ws.Range("AJ7").Value = 10971
$ws.Range("AK7").Value = 795000
$ws.Range("AL7").Value = 525424
$ws.Range("AM7").Value = 0.14
$ws.Range("AP7").Value = 11103
$ws.Range("AQ7").Value = 345
$ws.Range("AR7").Value = 5486
$ws.Range("AS7").Value = 141845
$ws.Range("AT7").Value = 0.2
$ws.Range("AV7").Value = 1
$ws.Range("AW7").Value = 5490
$ws.Range("AX7").Value = 1.68
$ws.Range("AZ7").Value = 525424
$ws.Range("BA7").Value = 4
$ws.Range("BB7").Value = 17
$ws.Range("BC7").Value = 16
$ws.Range("BD7").Value = 332
$ws.Range("BE7").Value = 5424

# Row 8
$ws.Range("Q8").Value = 8425
$ws.Range("R8").Value = 0
$ws.Range("Z8").Value = "8425"
$ws.Range("AA8").Value = "457 x 191 x 98"
$ws.Range("AB8").Value = 345
$ws.Range("AC8").Value = 8625
$ws.Range("AD8").Value = 7551
$ws.Range("AE8").Value = 1.07
$ws.Range("AG8").Value = 0.62
$ws.Range("AH8").Value = 5326
$ws.Range("AI8").Value = 345
$ws.Range("AJ8").Value = 8625
$ws.Range("AK8").Value = 625000
$ws.Range("AL8").Value = 119766
$ws.Range("AP8").Value = 8494
$ws.Range("AQ8").Value = 345
$ws.Range("AR8").Value = 4312
$ws.Range("AS8").Value = 76104
$ws.Range("AT8").Value = 0.24
$ws.Range("AV8").Value = 0.99
$ws.Range("AW8").Value = 4276
$ws.Range("AX8").Value = 1.42
$ws.Range("AZ8").Value = 103267
$ws.Range("BA8").Value = 3
$ws.Range("BB8").Value = 25
$ws.Range("BC8").Value = 47
$ws.Range("BD8").Value = 620
$ws.Range("BE8").Value = 4274

# Row 9
$ws.Range("Q9").Value = 8936
$ws.Range("R9").Value = 0
$ws.Range("Z9").Value = "8936"
$ws.Range("AA9").Value = "457 x 191 x 106"
$ws.Range("AB9").Value = 345
$ws.Range("AC9").Value = 9315
$ws.Range("AD9").Value = 32319
$ws.Range("AE9").Value = 0.54
$ws.Range("AH9").Value = 8499
$ws.Range("AI9").Value = 345
$ws.Range("AJ9").Value = 9315
$ws.Range("AK9").Value = 675000
$ws.Range("AL9").Value = 446115
$ws.Range("AM9").Value = 0.14
$ws.Range("AP9").Value = 9427
$ws.Range("AQ9").Value = 345
$ws.Range("AR9").Value = 4658
$ws.Range("AS9").Value = 81285
$ws.Range("AT9").Value = 0.24
$ws.Range("AW9").Value = 4617
$ws.Range("AX9").Value = 1.49
$ws.Range("AZ9").Value = 446115
$ws.Range("BA9").Value = 4
$ws.Range("BB9").Value = 17
$ws.Range("BC9").Value = 16
$ws.Range("BD9").Value = 277
$ws.Range("BE9").Value = 4496

# Row 10
$ws.Range("Q10").Value = 4953
$ws.Range("R10").Value = 0
$ws.Range("Z10").Value = "4953"
$ws.Range("AA10").Value = "457 x 152 x 60"
$ws.Range("AC10").Value = 5410
$ws.Range("AD10").Value = 4213
$ws.Range("AE10").Value = 1.13
$ws.Range("AG10").Value = 0.57
$ws.Range("AH10").Value = 3103
$ws.Range("AJ10").Value = 5410
$ws.Range("AK10").Value = 381000
$ws.Range("AL10").Value = 73009
$ws.Range("AP10").Value = 5323
$ws.Range("AR10").Value = 2705
$ws.Range("AS10").Value = 25746
$ws.Range("AT10").Value = 0.32
$ws.Range("AV10").Value = 0.97
$ws.Range("AW10").Value = 2629
$ws.Range("AX10").Value = 1.04
$ws.Range("AZ10").Value = 62952
$ws.Range("BA10").Value = 3
$ws.Range("BB10").Value = 25
$ws.Range("BC10").Value = 47
$ws.Range("BD10").Value = 385
$ws.Range("BE10").Value = 2515

# Row 11
$ws.Range("Q11").Value = 5254
$ws.Range("R11").Value = 0
$ws.Range("Z11").Value = "5254"
$ws.Range("AA11").Value = "533 x 165 x 66"
$ws.Range("AC11").Value = 5943
$ws.Range("AD11").Value = 23132
$ws.Range("AE11").Value = 0.51
$ws.Range("AG11").Value = 0.92
$ws.Range("AH11").Value = 5480
$ws.Range("AJ11").Value = 5943
$ws.Range("AK11").Value = 418500
$ws.Range("AL11").Value = 276591
$ws.Range("AP11").Value = 6012
$ws.Range("AR11").Value = 2971
$ws.Range("AS11").Value = 27818
$ws.Range("AT11").Value = 0.33
$ws.Range("AV11").Value = 0.97
$ws.Range("AW11").Value = 2886
$ws.Range("AX11").Value = 1.1
$ws.Range("AZ11").Value = 276591
$ws.Range("BA11").Value = 3
$ws.Range("BB11").Value = 13
$ws.Range("BC11").Value = 12
$ws.Range("BD11").Value = 165
$ws.Range("BE11").Value = 2644
